# Auto-generated: update market-price snapshot cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 78.833336
$ws.Range("I2").Value = 78.833336
$ws.Range("K2").Value = 78.833336
$ws.Range("M2").Value = 34.166664
$ws.Range("H8").Value = 110.166664
$ws.Range("I8").Value = 142.25
$ws.Range("J8").Value = 46
$ws.Range("K8").Value = 426.75
$ws.Range("L8").Value = 138
$ws.Range("M8").Value = -287.75
$ws.Range("N8").Value = -416
$ws.Range("H54").Value = 29950
$ws.Range("J54").Value = 29950
$ws.Range("L54").Value = 29950
$ws.Range("N54").Value = -30922
$ws.Range("H76").Value = 3151.182
$ws.Range("I76").Value = 2941.111
$ws.Range("J76").Value = 4096.5
$ws.Range("K76").Value = 2941.111
$ws.Range("L76").Value = 4096.5
$ws.Range("M76").Value = -2626.111
$ws.Range("N76").Value = -4726.5
$ws.Range("H79").Value = 3151.182
$ws.Range("I79").Value = 2941.111
$ws.Range("J79").Value = 4096.5
$ws.Range("K79").Value = 2941.111
$ws.Range("L79").Value = 4096.5
$ws.Range("M79").Value = -1849.111
$ws.Range("N79").Value = -6280.5
$ws.Range("H80").Value = 58830724
$ws.Range("I80").Value = 142857380
$ws.Range("J80").Value = 12064.2
$ws.Range("K80").Value = 428572140
$ws.Range("L80").Value = 36192.60000000001
$ws.Range("M80").Value = -428571142
$ws.Range("N80").Value = -38188.60000000001
$ws.Range("H81").Value = 117000
$ws.Range("J81").Value = 117000
$ws.Range("L81").Value = 117000
$ws.Range("N81").Value = -118996
$ws.Range("H83").Value = 58830724
$ws.Range("I83").Value = 142857380
$ws.Range("J83").Value = 12064.2
$ws.Range("K83").Value = 1285716420
$ws.Range("L83").Value = 108577.8
$ws.Range("M83").Value = -1285711428
$ws.Range("N83").Value = -118561.8
$ws.Range("H84").Value = 117000
$ws.Range("J84").Value = 117000
$ws.Range("L84").Value = 351000
$ws.Range("N84").Value = -360984
$ws.Range("H86").Value = 125004570
$ws.Range("J86").Value = 111115576
$ws.Range("L86").Value = 111115576
$ws.Range("N86").Value = -111117822
$ws.Range("H89").Value = 125004570
$ws.Range("J89").Value = 111115576
$ws.Range("L89").Value = 555577880
$ws.Range("N89").Value = -555589112
$ws.Range("H131").Value = 8464149
$ws.Range("I131").Value = 14286613
$ws.Range("J131").Value = 1671275
$ws.Range("K131").Value = 42859839
$ws.Range("L131").Value = 5013825
$ws.Range("M131").Value = -42854799
$ws.Range("N131").Value = -5023905
$ws.Range("H132").Value = 2998.6853
$ws.Range("I132").Value = 3019.8086
$ws.Range("K132").Value = 9059.425799999999
$ws.Range("M132").Value = -6529.425799999999
$ws.Range("H137").Value = 3028.5334
$ws.Range("I137").Value = 2154.48
$ws.Range("J137").Value = 7398.8
$ws.Range("K137").Value = 6463.440000000001
$ws.Range("L137").Value = 22196.4
$ws.Range("M137").Value = -3913.440000000001
$ws.Range("N137").Value = -27296.4
$ws.Range("H138").Value = 2669.7903
$ws.Range("J138").Value = 3282.8057
$ws.Range("L138").Value = 9848.417099999999
$ws.Range("N138").Value = -20128.4171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8174172
$ws.Range("I32").Value = 4505672
$ws.Range("J32").Value = 17869492
$ws.Range("K32").Value = 4505672
$ws.Range("L32").Value = 17869492
$ws.Range("M32").Value = -4505385
$ws.Range("N32").Value = -17870066
$ws.Range("H45").Value = 2623.6924
$ws.Range("I45").Value = 2209.1428
$ws.Range("K45").Value = 2209.1428
$ws.Range("M45").Value = -1832.1428
$ws.Range("H110").Value = 897
$ws.Range("I110").Value = 435.3889
$ws.Range("K110").Value = 435.3889
$ws.Range("M110").Value = 1609.6111
$ws.Range("H118").Value = 151994.5
$ws.Range("J118").Value = 151994.5
$ws.Range("L118").Value = 151994.5
$ws.Range("N118").Value = -155308.5
$ws.Range("H122").Value = 3234.8796
$ws.Range("I122").Value = 2318.05
$ws.Range("J122").Value = 5626.609
$ws.Range("K122").Value = 6954.150000000001
$ws.Range("L122").Value = 16879.827
$ws.Range("M122").Value = -4504.150000000001
$ws.Range("N122").Value = -21779.827

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2801.7058
$ws.Range("I99").Value = 2792.6428
$ws.Range("J99").Value = 2844
$ws.Range("K99").Value = 2792.6428
$ws.Range("L99").Value = 2844
$ws.Range("M99").Value = -1294.6428
$ws.Range("N99").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2402.2432
$ws.Range("I58").Value = 1858.8518
$ws.Range("J58").Value = 3869.4
$ws.Range("K58").Value = 1858.8518
$ws.Range("L58").Value = 3869.4
$ws.Range("M58").Value = -1655.8518
$ws.Range("N58").Value = -4275.4
$ws.Range("H99").Value = 2011.75
$ws.Range("I99").Value = 2011.75
$ws.Range("K99").Value = 2011.75
$ws.Range("M99").Value = -513.75
$ws.Range("H126").Value = 2011.75
$ws.Range("I126").Value = 2011.75
$ws.Range("K126").Value = 6035.25
$ws.Range("M126").Value = -3565.25
$ws.Range("H134").Value = 3738.8262
$ws.Range("I134").Value = 2531.4375
$ws.Range("K134").Value = 7594.3125
$ws.Range("M134").Value = -5059.3125
$ws.Range("H136").Value = 2402.2432
$ws.Range("I136").Value = 1858.8518
$ws.Range("J136").Value = 3869.4
$ws.Range("K136").Value = 5576.555399999999
$ws.Range("L136").Value = 11608.2
$ws.Range("M136").Value = -3026.555399999999
$ws.Range("N136").Value = -16708.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 671.2857
$ws.Range("J51").Value = 740.2
$ws.Range("L51").Value = 2220.6
$ws.Range("N51").Value = -3140.6
$ws.Range("H61").Value = 539.25
$ws.Range("J61").Value = 539.25
$ws.Range("L61").Value = 1617.75
$ws.Range("N61").Value = -2047.75
$ws.Range("H68").Value = 1531.8889
$ws.Range("J68").Value = 1677.5
$ws.Range("L68").Value = 5032.5
$ws.Range("N68").Value = -6654.5
$ws.Range("H71").Value = 1531.8889
$ws.Range("J71").Value = 1677.5
$ws.Range("L71").Value = 15097.5
$ws.Range("N71").Value = -23209.5
$ws.Range("H112").Value = 3159.2
$ws.Range("H113").Value = 1269.2778
$ws.Range("J113").Value = 1348.625
$ws.Range("L113").Value = 4045.875
$ws.Range("N113").Value = -8385.875
$ws.Range("H122").Value = 739.0952
$ws.Range("I122").Value = 545.5
$ws.Range("J122").Value = 816.5333000000001
$ws.Range("K122").Value = 4909.5
$ws.Range("L122").Value = 7348.7997
$ws.Range("M122").Value = -2459.5
$ws.Range("N122").Value = -12248.7997
$ws.Range("H123").Value = 2724.75
$ws.Range("I123").Value = 2450
$ws.Range("K123").Value = 7350
$ws.Range("M123").Value = -4900
$ws.Range("H129").Value = 1612
$ws.Range("J129").Value = 2010.3334
$ws.Range("L129").Value = 6031.0002
$ws.Range("N129").Value = -16031.0002
$ws.Range("H131").Value = 1547.1111
$ws.Range("I131").Value = 1211.3846
$ws.Range("J131").Value = 1736.8695
$ws.Range("K131").Value = 3634.1538
$ws.Range("L131").Value = 5210.6085
$ws.Range("M131").Value = 1405.8462
$ws.Range("N131").Value = -15290.6085

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1853.5
$ws.Range("I122").Value = 1677.0667
$ws.Range("K122").Value = 5031.2001
$ws.Range("M122").Value = -2581.2001
$ws.Range("H132").Value = 7166.3335
$ws.Range("I132").Value = 7500
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -19970
$ws.Range("N132").Value = -26058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14577.414
$ws.Range("I7").Value = 14154.523
$ws.Range("J7").Value = 15687.5
$ws.Range("K7").Value = 14154.523
$ws.Range("L7").Value = 15687.5
$ws.Range("M7").Value = -14042.523
$ws.Range("N7").Value = -15911.5
$ws.Range("H40").Value = 5857.3716
$ws.Range("I40").Value = 5583.7
$ws.Range("K40").Value = 5583.7
$ws.Range("M40").Value = -5447.7
$ws.Range("H46").Value = 2595.875
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 2895.2856
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 2895.2856
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -3271.2856
$ws.Range("H126").Value = 14577.414
$ws.Range("I126").Value = 14154.523
$ws.Range("J126").Value = 15687.5
$ws.Range("K126").Value = 42463.569
$ws.Range("L126").Value = 47062.5
$ws.Range("M126").Value = -39993.569
$ws.Range("N126").Value = -52002.5
$ws.Range("H141").Value = 363333
$ws.Range("J141").Value = 524999.5
$ws.Range("L141").Value = 524999.5
$ws.Range("N141").Value = -535359.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2589.9048
$ws.Range("I122").Value = 2011.7059
$ws.Range("K122").Value = 6035.1177
$ws.Range("M122").Value = -3585.1177
$ws.Range("H123").Value = 110994
$ws.Range("J123").Value = 110994
$ws.Range("L123").Value = 110994
$ws.Range("N123").Value = -120794
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820
$ws.Range("H125").Value = 99874.5
$ws.Range("J125").Value = 99874.5
$ws.Range("L125").Value = 99874.5
$ws.Range("N125").Value = -109714.5
$ws.Range("H128").Value = 100000000
$ws.Range("J128").Value = 100000000
$ws.Range("L128").Value = 100000000
$ws.Range("N128").Value = -100009960
$ws.Range("H132").Value = 4843.2915
$ws.Range("I132").Value = 4176.7856
$ws.Range("J132").Value = 5776.4
$ws.Range("K132").Value = 12530.3568
$ws.Range("L132").Value = 17329.2
$ws.Range("M132").Value = -10000.3568
$ws.Range("N132").Value = -22389.2

Write-Output "Applied 250 cell updates"